$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data record (row 2) being added to the collections index.
$ws.Range("A2").Value = "MCH114"
$ws.Range("C2").Value = "CORRESPONDENCE: TRUDIE GELB, PUBLICATIONS: LOOKING FORWARD, WORKERS IN CHAINS, ECONOMIC CRISIS IN SA, WORKERS STRUGGLE FOR FREEDOM, ROOTING DEMOCRACY ON AFRICAN SOIL, CONSTITUTIONAL STRUCTURES, LOCAL GOVERNMENT & PLANNING 1992 SECOND ANNUAL REPORT"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1992"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Row 2 uses the worksheet's secondary body style: 10pt Calibri, theme
# text color, no fill - matching the rest of the data rows.
foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

$ws.Range("B6").Select()
